$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-0.8333010676406252, 1.032645219681061, 0.6420169025845999, 2.6108000046005815)"
$ws.Range("C2").Value = "NIG(0.9348143054394998, 0.6573978112080057, 4.431484423667587, 5.057080539629593)"
$ws.Range("D2").Value = "NCT(2.925058423808954, 2.192370767486075, -1.933079655429136, 2.217432491158548)"
$ws.Range("E2").Value = "JSU(-1.270267578414953, 1.2257364308782934, 2.6356580812197556, 5.611613241333833)"
